$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 570 (before the old row 571),
# shifting the existing "2026/12/29 ..." block and everything below it
# down by two rows.
$ws.Rows.Item(571).Resize(2).Insert()

# Populate the two newly inserted rows (571, 572) with the new
# "2026/01/07" / "水" entries.
#
# Column A holds dates stored as plain text (no explicit number format in
# this workbook), so a leading apostrophe is used to force text entry -
# otherwise Excel would silently reinterpret the date-shaped string as a
# real date serial number. ClearFormats() afterwards drops the implicit
# "quote prefix" formatting that the forced text-entry leaves behind, so
# the cell ends up with the same default (un-styled) look as every other
# date cell in the column.
$ws.Cells.Item(571, 1).Value = "'2026/01/07"
$ws.Cells.Item(571, 1).ClearFormats()
$ws.Cells.Item(571, 2).Value = "水"
$ws.Cells.Item(571, 3).Value = 7
$ws.Cells.Item(571, 4).Value = 186

$ws.Cells.Item(572, 1).Value = "'2026/01/07"
$ws.Cells.Item(572, 1).ClearFormats()
$ws.Cells.Item(572, 2).Value = "水"
$ws.Cells.Item(572, 3).Value = 10
$ws.Cells.Item(572, 4).Value = 189
